# Applies the "Error Calculations and Plots" edit:
# - Toggle several missing/present values in column E (rows 6,8,12,14,17,18,19,20,23)
# - Remove rows for "RM 232" and "SC 92" (data points dropped), shifting subsequent rows up
# - Update a few column C / E values in some of the shifted rows
# - Dimension shrinks from A1:F35 to A1:F33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) In-place toggles on column E (before the row deletions, since these rows are above row 26)
$ws.Cells.Item(6, 5).Value = -5.7          # E6: was blank -> -5.7
$ws.Cells.Item(8, 5).ClearContents()       # E8: was -6.6 -> blank
$ws.Cells.Item(12, 5).Value = -5.3         # E12: was blank -> -5.3
$ws.Cells.Item(14, 5).ClearContents()      # E14: was -5.4 -> blank
$ws.Cells.Item(17, 5).Value = -7.3         # E17: was blank -> -7.3
$ws.Cells.Item(18, 5).Value = -8.5         # E18: was blank -> -8.5
$ws.Cells.Item(19, 5).ClearContents()      # E19: was -6.5 -> blank
$ws.Cells.Item(20, 5).ClearContents()      # E20: was -7.2 -> blank
$ws.Cells.Item(23, 5).Value = -7           # E23: was blank -> -7

# 2) Remove the "RM 232" row (originally row 26) entirely
$ws.Rows.Item(26).Delete()

# 3) Remove the "SC 92" row. After the previous delete it now sits at row 27
$ws.Rows.Item(27).Delete()

# 4) Tweak a few values on rows that shifted up into their new positions
# "SC 101" is now row 27
$ws.Cells.Item(27, 3).Value = 10           # C27: was blank -> 10
$ws.Cells.Item(27, 5).ClearContents()      # E27: was -10 -> blank

# "SC 119" is now row 29
$ws.Cells.Item(29, 3).ClearContents()      # C29: was 11.2 -> blank

# "SC 193" is now row 32
$ws.Cells.Item(32, 3).ClearContents()      # C32: was 10.5 -> blank
